$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "             DATE"
$ws.Range("B1").Value = "          TASK NAME"
$ws.Range("C1").Value = "    LEARNING FROM TASK"
$ws.Range("D1").Value = "       EXTRA STUDY"

# Row 3
$ws.Range("A3").Value = "      05-03-2024"
$ws.Range("B3").Value = " Login Page using HTML CSS"
$ws.Range("C3").Value = " Basic HTML page creation"
$ws.Range("D3").Value = " I started leaning about Typescript"

# Row 4
$ws.Range("A4").Value = "      06-03-2024"
$ws.Range("B4").Value = "Create Registration Page "
$ws.Range("C4").Value = " Connecting different HTML pages"
$ws.Range("D4").Value = "Continued Typescript"

# Row 5
$ws.Range("A5").Value = "      07-03-2024"
$ws.Range("B5").Value = "Forget Module HTML CSS"
$ws.Range("C5").Value = "Connected on more Page"
$ws.Range("D5").Value = "Studies About Nextjs"

# Row 6
$ws.Range("A6").Value = "      08-03-2024"
$ws.Range("B6").Value = "Forget Module OTP VERIFY"
$ws.Range("C6").Value = "OTP verification using Node"
$ws.Range("D6").Value = "Learned About Node Mailer"

# Row 9
$ws.Range("A9").Value = "      11-03-2024"
$ws.Range("B9").Value = "Landing Page 5 Star Hotel"
$ws.Range("C9").Value = "swipe js, animista-css library explore"

# Row 10
$ws.Range("A10").Value = "      12-03-2024"
$ws.Range("B10").Value = "Hotel Management Backend"
$ws.Range("C10").Value = "modules,controllers,routes in Node JS"

# Row 11
$ws.Range("A11").Value = "      13-03-2024"
$ws.Range("B11").Value = "Hotel Management Backend"
$ws.Range("C11").Value = "Postment and API testing"

# Column widths (closest achievable values given engine's internal width
# quantization; targets are 13.77734375 / 24.5546875 / 39.77734375 / 31.5546875)
$ws.Columns.Item(1).ColumnWidth = 13.0
$ws.Columns.Item(2).ColumnWidth = 23.666666666666668
$ws.Columns.Item(3).ColumnWidth = 39.0
$ws.Columns.Item(4).ColumnWidth = 30.666666666666668

# Header/Footer
$ws.PageSetup.CenterHeader = "DAILY LEARNING"
$ws.PageSetup.Orientation = 1

# Selection similar to target
$null = $ws.Range("C13").Select()

Write-Output "Applied daily-learning edits"
